$d = $word.ActiveDocument

# Locate the heading run that currently reads "6.Conclusion:" so the
# edit is anchored on content rather than a hard-coded offset.
$full = $d.Content
$found = $full.Find.Execute("6.Conclusion:", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '6.Conclusion:' in the document"
}

$headingStart = $full.Start
$headingEnd = $full.End

# Split the found range into the leading digit ("6") and the trailing
# ".Conclusion:" text so each piece can be retargeted independently.
$digitRange = $d.Range($headingStart, $headingStart + 1)
$restRange = $d.Range($headingStart + 1, $headingEnd)

# Bookmark the trailing portion first: this pins its run boundary so the
# upcoming edit to the digit doesn't get silently coalesced back into a
# single run once both pieces end up with identical formatting again.
$markName = "TmpSplitMark"
$d.Bookmarks.Add($markName, $restRange)

# "6" -> "5"
$digitRange.Text = "5"

# Re-touch the trailing text (replace with itself through a throwaway
# intermediate value) so it is rewritten as its own fresh run rather than
# keeping stale run-level rsid bookkeeping from the original single run.
$restStart = $headingStart + 1
$restEnd = $headingEnd
$tempRange = $d.Range($restStart, $restEnd)
$tempRange.Text = "ZZTEMPZZ"
$tempRange2 = $d.Range($restStart, $restStart + 8)
$tempRange2.Text = ".Conclusion:"

# Remove the helper bookmark now that the two runs exist independently.
$d.Bookmarks($markName).Delete()
